# Update mapa_interactivo_INCO.xlsx ("INCO" sheet):
# Two report rows were removed from the source data export:
#   - Row with Caso -171 (OLAVARRIA /ALT/ 531)           -> originally sheet row 19
#   - Row with Caso -256 (NECOCHEA /ALT/ 1279)            -> originally sheet row 45
# Removing full rows shifts every row below up by one, which is exactly what
# the target workbook shows (dimension goes from A1:N50 to A1:N48).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

# Delete the lower row first so the row index of the first deletion (19)
# doesn't move before we get to use it.
$ws.Rows.Item(45).EntireRow.Delete()
$ws.Rows.Item(19).EntireRow.Delete()
